$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "29.365.14"
Set-TextCell $ws.Range("E2") "  +0.63%  "

Set-TextCell $ws.Range("D3") "1.874.88"
Set-TextCell $ws.Range("E3") "  +0.91%  "

Set-TextCell $ws.Range("E4") "  -0.03%  "

Set-TextCell $ws.Range("D5") "0.7126"
Set-TextCell $ws.Range("E5") "  +0.38%  "

Set-TextCell $ws.Range("D6") "241.85"
Set-TextCell $ws.Range("E6") "  +0.72%  "

Set-TextCell $ws.Range("E7") "  -0.06%  "

Set-TextCell $ws.Range("D8") "0.07799"
Set-TextCell $ws.Range("E8") "  +2.25%  "

Set-TextCell $ws.Range("D9") "0.3114"
Set-TextCell $ws.Range("E9") "  +1.54%  "

Set-TextCell $ws.Range("D10") "25.23"
Set-TextCell $ws.Range("E10") "  +2.38%  "

Set-TextCell $ws.Range("D11") "0.08443"
Set-TextCell $ws.Range("E11") "  +0.99%  "

Set-TextCell $ws.Range("D12") "1.873.66"
Set-TextCell $ws.Range("E12") "  +1.09%  "

Set-TextCell $ws.Range("D13") "5.245"
Set-TextCell $ws.Range("E13") "  +1.40%  "

Set-TextCell $ws.Range("D14") "0.7141"
Set-TextCell $ws.Range("E14") "  +1.07%  "

Set-TextCell $ws.Range("D15") "91.16"
Set-TextCell $ws.Range("E15") "  +0.08%  "

Set-TextCell $ws.Range("D16") "29.367.97"
Set-TextCell $ws.Range("E16") "  +0.61%  "

Set-TextCell $ws.Range("D17") "6.079"
Set-TextCell $ws.Range("E17") "  +2.58%  "

Set-TextCell $ws.Range("D18") "0.000008248"
Set-TextCell $ws.Range("E18") "  +5.78%  "

Set-TextCell $ws.Range("D19") "241.02"
Set-TextCell $ws.Range("E19") "  -0.44%  "

Set-TextCell $ws.Range("D20") "13.26"
Set-TextCell $ws.Range("E20") "  +1.32%  "

Set-TextCell $ws.Range("D21") "2.122.65"
Set-TextCell $ws.Range("E21") "  +0.25%  "

Set-TextCell $ws.Range("D22") "1.000"
Set-TextCell $ws.Range("E22") "  -0.05%  "

Set-TextCell $ws.Range("D23") "7.786"
Set-TextCell $ws.Range("E23") "  -0.54%  "

Set-TextCell $ws.Range("E24") "  -0.04%  "

Set-TextCell $ws.Range("D25") "0.1595"
Set-TextCell $ws.Range("E25") "  +0.60%  "

Set-TextCell $ws.Range("D26") "163.14"
Set-TextCell $ws.Range("E26") "  +0.47%  "

Set-TextCell $ws.Range("E27") "  +2.29%  "

Set-TextCell $ws.Range("E28") "  +0.74%  "

Set-TextCell $ws.Range("E29") "  +1.10%  "

Set-TextCell $ws.Range("D30") "4.425"
Set-TextCell $ws.Range("E30") "  +1.04%  "

Set-TextCell $ws.Range("D31") "4.329"

Set-TextCell $ws.Range("D32") "1.288"
Set-TextCell $ws.Range("E32") "  -3.13%  "

Set-TextCell $ws.Range("D33") "0.05304"
Set-TextCell $ws.Range("E33") "  +4.02%  "

Set-TextCell $ws.Range("D34") "1.939"
Set-TextCell $ws.Range("E34") "  +1.39%  "

Set-TextCell $ws.Range("D35") "1.179"
Set-TextCell $ws.Range("E35") "  +1.64%  "

Set-TextCell $ws.Range("D36") "0.7447"
Set-TextCell $ws.Range("E36") "  -6.82%  "

Set-TextCell $ws.Range("E37") "  +0.51%  "

Set-TextCell $ws.Range("D38") "0.01870"
Set-TextCell $ws.Range("E38") "  +1.62%  "

Set-TextCell $ws.Range("D39") "1.228.81"
Set-TextCell $ws.Range("E39") "  +4.39%  "

Set-TextCell $ws.Range("D40") "2.731"

Set-TextCell $ws.Range("D41") "6.522"
Set-TextCell $ws.Range("E41") "  +5.94%  "

Set-TextCell $ws.Range("D42") "111.15"
Set-TextCell $ws.Range("E42") "  +9.42%  "

Set-TextCell $ws.Range("D43") "0.8923"
Set-TextCell $ws.Range("E43") "  +0.35%  "

Set-TextCell $ws.Range("D44") "72.98"
Set-TextCell $ws.Range("E44") "  +0.81%  "

Set-TextCell $ws.Range("E45") "  +0.02%  "

Set-TextCell $ws.Range("D46") "2.020.48"
Set-TextCell $ws.Range("E46") "  +0.30%  "

Set-TextCell $ws.Range("D47") "1.815"
Set-TextCell $ws.Range("E47") "  +2.41%  "

Set-TextCell $ws.Range("B48") "Mantle"
Set-TextCell $ws.Range("C48") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws.Range("D48") "0.5214"
Set-TextCell $ws.Range("E48") "  +0.74%  "

Set-TextCell $ws.Range("B49") "BabyDogeCoin"
Set-TextCell $ws.Range("C49") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell $ws.Range("D49") "0.00000000123"
Set-TextCell $ws.Range("E49") "  +2.51%  "

Set-TextCell $ws.Range("D50") "9.433"
Set-TextCell $ws.Range("E50") "  +2.30%  "

Set-TextCell $ws.Range("D51") "0.4324"
Set-TextCell $ws.Range("E51") "  +1.88%  "
